# The presentation currently uses the "Integral" (Red Violet) theme for its
# slides/slide master, while the plain default "Office Theme" palette sits
# unused on the notes master. The author switched the deck's design back to
# the default Office colour theme, so every themed colour slot (the 12
# modern theme colours: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) needs to be
# set back to the stock Office palette.
#
# PowerPoint exposes those 12 slots via Slide.ThemeColorScheme(1..12), in
# the standard msoThemeColorDark1..msoThemeColorFollowedHyperlink order.
# Writing through any slide updates the single shared theme part backing
# the whole deck's design, so slide 1 is as good an anchor as any.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> RGB() value (R + G*256 + B*65536) for the stock "Office" theme.
$officeColors = @{
    1  = 0          # dk1      #000000
    2  = 16777215   # lt1      #FFFFFF
    3  = 6968388    # dk2      #44546A
    4  = 15132391   # lt2      #E7E6E6
    5  = 13998939   # accent1  #5B9BD5
    6  = 3243501     # accent2  #ED7D31
    7  = 10855845   # accent3  #A5A5A5
    8  = 49407       # accent4  #FFC000
    9  = 12874308   # accent5  #4472C4
    10 = 4697456     # accent6  #70AD47
    11 = 12673797   # hlink    #0563C1
    12 = 7491477     # folHlink #954F72
}

foreach ($idx in 1..12) {
    $tcs.Item($idx).RGB = $officeColors[$idx]
}

# Best-effort: restore the friendly design label too ("Integral" ->
# "Office Theme"). Harmless if the host doesn't expose this as writable.
try { $p.SlideMaster.Theme.Name = "Office Theme" } catch {}
